# Insert a new row at row 79 (this shifts existing rows 79-86 down to 80-87,
# preserving their content/formatting), then fill the new row 79 with the
# new weekly price entry for Ajo (garlic) at Vega Monumental Concepcion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(79).Insert()

$newRow = 79
$ws.Cells.Item($newRow, 1).Value = 11
$ws.Cells.Item($newRow, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item($newRow, 3).Value = "Bíobío"
$ws.Cells.Item($newRow, 4).Value = 44449
$ws.Cells.Item($newRow, 4).NumberFormat = $ws.Cells.Item($newRow + 1, 4).NumberFormat
$ws.Cells.Item($newRow, 5).Value = 8
$ws.Cells.Item($newRow, 6).Value = 100112003
$ws.Cells.Item($newRow, 7).Value = "Ajo"
$ws.Cells.Item($newRow, 8).Value = "Chino"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 400
$ws.Cells.Item($newRow, 11).Value = 15000
$ws.Cells.Item($newRow, 12).Value = 15500
$ws.Cells.Item($newRow, 13).Value = 15250
$ws.Cells.Item($newRow, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item($newRow, 15).Value = "China"
$ws.Cells.Item($newRow, 16).Value = 1525
$ws.Cells.Item($newRow, 17).Value = 10
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
